$wb = $excel.ActiveWorkbook

# --- Update the "conversion" text on sheet Hoja1 (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.03 = 7755.49 pesos`n✅ 7755.49 pesos = 2.02 = 949.81 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate cells on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 492.998
$wsTasas.Range("O10").Value = 3823.44
$wsTasas.Range("N12").Value = 3845.88
$wsTasas.Range("O12").Value = 471.001
